# Batch cohort generation - WIP
# Adds "Dist" (distribution) and "CV" columns to the parameter table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the Dist column
$ws.Range("G1").Value = "Dist"

# Fill Dist values first so the "norm" shared string is registered
# before the "CV" header string.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 7).Value = "norm"
}

# New header cell for the CV column
$ws.Range("H1").Value = "CV"

# Fill CV values for each parameter row (rows 2-9)
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 8).Value = 0.05
}

# Update active selection to match the authored state
$ws.Range("H12").Select()
